$wb = $excel.ActiveWorkbook

# ---- mmWave(HR) ----
$ws = $wb.Worksheets.Item("mmWave(HR)")
$srcRow = 73
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A74:F74"))
$ws.Cells.Item(74, 2).Value = "18:32:59"
$ws.Cells.Item(74, 3).Value = "18:00"
$ws.Cells.Item(74, 4).Value = "Bedroom"
$ws.Cells.Item(74, 5).Value = 2
$ws.Cells.Item(74, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A75:F75"))
$ws.Cells.Item(75, 2).Value = "18:33:02"
$ws.Cells.Item(75, 3).Value = "18:00"
$ws.Cells.Item(75, 4).Value = "Bedroom"
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A76:F76"))
$ws.Cells.Item(76, 2).Value = "18:33:05"
$ws.Cells.Item(76, 3).Value = "18:00"
$ws.Cells.Item(76, 4).Value = "Bedroom"
$ws.Cells.Item(76, 5).Value = 2
$ws.Cells.Item(76, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A77:F77"))
$ws.Cells.Item(77, 2).Value = "18:33:10"
$ws.Cells.Item(77, 3).Value = "18:00"
$ws.Cells.Item(77, 4).Value = "Bedroom"
$ws.Cells.Item(77, 5).Value = 1
$ws.Cells.Item(77, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A78:F78"))
$ws.Cells.Item(78, 2).Value = "18:33:15"
$ws.Cells.Item(78, 3).Value = "18:00"
$ws.Cells.Item(78, 4).Value = "Bedroom"
$ws.Cells.Item(78, 5).Value = 2
$ws.Cells.Item(78, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A79:F79"))
$ws.Cells.Item(79, 2).Value = "18:33:18"
$ws.Cells.Item(79, 3).Value = "18:00"
$ws.Cells.Item(79, 4).Value = "Bedroom"
$ws.Cells.Item(79, 5).Value = 27
$ws.Cells.Item(79, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A80:F80"))
$ws.Cells.Item(80, 2).Value = "18:33:21"
$ws.Cells.Item(80, 3).Value = "18:00"
$ws.Cells.Item(80, 4).Value = "Bedroom"
$ws.Cells.Item(80, 5).Value = 2
$ws.Cells.Item(80, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A81:F81"))
$ws.Cells.Item(81, 2).Value = "18:33:24"
$ws.Cells.Item(81, 3).Value = "18:00"
$ws.Cells.Item(81, 4).Value = "Bedroom"
$ws.Cells.Item(81, 5).Value = 1
$ws.Cells.Item(81, 6).Value = "Occupied"

# ---- mmWave(BR) ----
$ws = $wb.Worksheets.Item("mmWave(BR)")
$srcRow = 73
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A74:F74"))
$ws.Cells.Item(74, 2).Value = "18:32:58"
$ws.Cells.Item(74, 3).Value = "18:00"
$ws.Cells.Item(74, 4).Value = "Bedroom"
$ws.Cells.Item(74, 5).Value = 50
$ws.Cells.Item(74, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A75:F75"))
$ws.Cells.Item(75, 2).Value = "18:33:01"
$ws.Cells.Item(75, 3).Value = "18:00"
$ws.Cells.Item(75, 4).Value = "Bedroom"
$ws.Cells.Item(75, 5).Value = 52
$ws.Cells.Item(75, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A76:F76"))
$ws.Cells.Item(76, 2).Value = "18:33:04"
$ws.Cells.Item(76, 3).Value = "18:00"
$ws.Cells.Item(76, 4).Value = "Bedroom"
$ws.Cells.Item(76, 5).Value = 50
$ws.Cells.Item(76, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A77:F77"))
$ws.Cells.Item(77, 2).Value = "18:33:09"
$ws.Cells.Item(77, 3).Value = "18:00"
$ws.Cells.Item(77, 4).Value = "Bedroom"
$ws.Cells.Item(77, 5).Value = 49
$ws.Cells.Item(77, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A78:F78"))
$ws.Cells.Item(78, 2).Value = "18:33:13"
$ws.Cells.Item(78, 3).Value = "18:00"
$ws.Cells.Item(78, 4).Value = "Bedroom"
$ws.Cells.Item(78, 5).Value = 50
$ws.Cells.Item(78, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A79:F79"))
$ws.Cells.Item(79, 2).Value = "18:33:17"
$ws.Cells.Item(79, 3).Value = "18:00"
$ws.Cells.Item(79, 4).Value = "Bedroom"
$ws.Cells.Item(79, 5).Value = 75
$ws.Cells.Item(79, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A80:F80"))
$ws.Cells.Item(80, 2).Value = "18:33:20"
$ws.Cells.Item(80, 3).Value = "18:00"
$ws.Cells.Item(80, 4).Value = "Bedroom"
$ws.Cells.Item(80, 5).Value = 50
$ws.Cells.Item(80, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A81:F81"))
$ws.Cells.Item(81, 2).Value = "18:33:23"
$ws.Cells.Item(81, 3).Value = "18:00"
$ws.Cells.Item(81, 4).Value = "Bedroom"
$ws.Cells.Item(81, 5).Value = 49
$ws.Cells.Item(81, 6).Value = "Occupied"

# ---- mmWave(InBed) ----
$ws = $wb.Worksheets.Item("mmWave(InBed)")
$srcRow = 73
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A74:F74"))
$ws.Cells.Item(74, 2).Value = "18:32:57"
$ws.Cells.Item(74, 3).Value = "18:00"
$ws.Cells.Item(74, 4).Value = "Bedroom"
$ws.Cells.Item(74, 5).Value = "In Bed"
$ws.Cells.Item(74, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A75:F75"))
$ws.Cells.Item(75, 2).Value = "18:33:00"
$ws.Cells.Item(75, 3).Value = "18:00"
$ws.Cells.Item(75, 4).Value = "Bedroom"
$ws.Cells.Item(75, 5).Value = "In Bed"
$ws.Cells.Item(75, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A76:F76"))
$ws.Cells.Item(76, 2).Value = "18:33:03"
$ws.Cells.Item(76, 3).Value = "18:00"
$ws.Cells.Item(76, 4).Value = "Bedroom"
$ws.Cells.Item(76, 5).Value = "In Bed"
$ws.Cells.Item(76, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A77:F77"))
$ws.Cells.Item(77, 2).Value = "18:33:08"
$ws.Cells.Item(77, 3).Value = "18:00"
$ws.Cells.Item(77, 4).Value = "Bedroom"
$ws.Cells.Item(77, 5).Value = "In Bed"
$ws.Cells.Item(77, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A78:F78"))
$ws.Cells.Item(78, 2).Value = "18:33:12"
$ws.Cells.Item(78, 3).Value = "18:00"
$ws.Cells.Item(78, 4).Value = "Bedroom"
$ws.Cells.Item(78, 5).Value = "In Bed"
$ws.Cells.Item(78, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A79:F79"))
$ws.Cells.Item(79, 2).Value = "18:33:16"
$ws.Cells.Item(79, 3).Value = "18:00"
$ws.Cells.Item(79, 4).Value = "Bedroom"
$ws.Cells.Item(79, 5).Value = "In Bed"
$ws.Cells.Item(79, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A80:F80"))
$ws.Cells.Item(80, 2).Value = "18:33:19"
$ws.Cells.Item(80, 3).Value = "18:00"
$ws.Cells.Item(80, 4).Value = "Bedroom"
$ws.Cells.Item(80, 5).Value = "In Bed"
$ws.Cells.Item(80, 6).Value = "Occupied"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A81:F81"))
$ws.Cells.Item(81, 2).Value = "18:33:22"
$ws.Cells.Item(81, 3).Value = "18:00"
$ws.Cells.Item(81, 4).Value = "Bedroom"
$ws.Cells.Item(81, 5).Value = "In Bed"
$ws.Cells.Item(81, 6).Value = "Occupied"

# ---- Proximity ----
$ws = $wb.Worksheets.Item("Proximity")
$srcRow = 23
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A24:F24"))
$ws.Cells.Item(24, 2).Value = "18:32:56"
$ws.Cells.Item(24, 3).Value = "18:00"
$ws.Cells.Item(24, 4).Value = "Living Room Main Door"
$ws.Cells.Item(24, 5).Value = "EXIT"
$ws.Cells.Item(24, 6).Value = "User EXITED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A25:F25"))
$ws.Cells.Item(25, 2).Value = "18:33:06"
$ws.Cells.Item(25, 3).Value = "18:00"
$ws.Cells.Item(25, 4).Value = "Living Room Main Door"
$ws.Cells.Item(25, 5).Value = "ENTER"
$ws.Cells.Item(25, 6).Value = "User ENTERED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A26:F26"))
$ws.Cells.Item(26, 2).Value = "18:33:11"
$ws.Cells.Item(26, 3).Value = "18:00"
$ws.Cells.Item(26, 4).Value = "Living Room Main Door"
$ws.Cells.Item(26, 5).Value = "EXIT"
$ws.Cells.Item(26, 6).Value = "User EXITED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A27:F27"))
$ws.Cells.Item(27, 2).Value = "18:33:25"
$ws.Cells.Item(27, 3).Value = "18:00"
$ws.Cells.Item(27, 4).Value = "Living Room Main Door"
$ws.Cells.Item(27, 5).Value = "ENTER"
$ws.Cells.Item(27, 6).Value = "User ENTERED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A28:F28"))
$ws.Cells.Item(28, 2).Value = "18:33:27"
$ws.Cells.Item(28, 3).Value = "18:00"
$ws.Cells.Item(28, 4).Value = "Living Room Main Door"
$ws.Cells.Item(28, 5).Value = "EXIT"
$ws.Cells.Item(28, 6).Value = "User EXITED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A29:F29"))
$ws.Cells.Item(29, 2).Value = "18:33:29"
$ws.Cells.Item(29, 3).Value = "18:00"
$ws.Cells.Item(29, 4).Value = "Living Room Main Door"
$ws.Cells.Item(29, 5).Value = "ENTER"
$ws.Cells.Item(29, 6).Value = "User ENTERED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A30:F30"))
$ws.Cells.Item(30, 2).Value = "18:33:43"
$ws.Cells.Item(30, 3).Value = "18:00"
$ws.Cells.Item(30, 4).Value = "Living Room Main Door"
$ws.Cells.Item(30, 5).Value = "EXIT"
$ws.Cells.Item(30, 6).Value = "User EXITED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A31:F31"))
$ws.Cells.Item(31, 2).Value = "18:33:46"
$ws.Cells.Item(31, 3).Value = "18:00"
$ws.Cells.Item(31, 4).Value = "Living Room Main Door"
$ws.Cells.Item(31, 5).Value = "ENTER"
$ws.Cells.Item(31, 6).Value = "User ENTERED Living Room Main Door"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A32:F32"))
$ws.Cells.Item(32, 2).Value = "18:33:51"
$ws.Cells.Item(32, 3).Value = "18:00"
$ws.Cells.Item(32, 4).Value = "Living Room Main Door"
$ws.Cells.Item(32, 5).Value = "EXIT"
$ws.Cells.Item(32, 6).Value = "User EXITED Living Room Main Door"

# ---- Camera ----
$ws = $wb.Worksheets.Item("Camera")
$srcRow = 14
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A15:F15"))
$ws.Cells.Item(15, 2).Value = "18:33:07"
$ws.Cells.Item(15, 3).Value = "18:00"
$ws.Cells.Item(15, 4).Value = "Living Room Main Door"
$ws.Cells.Item(15, 5).Value = "Image Captured"
$ws.Cells.Item(15, 6).Value = "Active"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A16:F16"))
$ws.Cells.Item(16, 2).Value = "18:33:26"
$ws.Cells.Item(16, 3).Value = "18:00"
$ws.Cells.Item(16, 4).Value = "Living Room Main Door"
$ws.Cells.Item(16, 5).Value = "Image Captured"
$ws.Cells.Item(16, 6).Value = "Active"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A17:F17"))
$ws.Cells.Item(17, 2).Value = "18:33:28"
$ws.Cells.Item(17, 3).Value = "18:00"
$ws.Cells.Item(17, 4).Value = "Living Room Main Door"
$ws.Cells.Item(17, 5).Value = "Image Captured"
$ws.Cells.Item(17, 6).Value = "Active"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A18:F18"))
$ws.Cells.Item(18, 2).Value = "18:33:30"
$ws.Cells.Item(18, 3).Value = "18:00"
$ws.Cells.Item(18, 4).Value = "Living Room Main Door"
$ws.Cells.Item(18, 5).Value = "Image Captured"
$ws.Cells.Item(18, 6).Value = "Active"
$ws.Range("A" + $srcRow + ":F" + $srcRow).Copy($ws.Range("A19:F19"))
$ws.Cells.Item(19, 2).Value = "18:33:50"
$ws.Cells.Item(19, 3).Value = "18:00"
$ws.Cells.Item(19, 4).Value = "Living Room Main Door"
$ws.Cells.Item(19, 5).Value = "Image Captured"
$ws.Cells.Item(19, 6).Value = "Active"

